$d = $word.ActiveDocument

# Bullet: "Closing 2 and 3 runs successively reduces (though not significant)
# ticket price and revenue." -> "...reduce (though not significant)..."
# ("reduces" -> "reduce").
$ok1 = $d.Content.Find.Execute("reduces (though not significant)", $true, $false, $false, $false, $false,
                                $true, 1, $false, "reduce (though not significant)", 2)
if (-not $ok1) {
    throw "Could not find target text for the 'reduces' -> 'reduce' correction."
}

# Bullet: "Closing 3, 4 or 5 runs has same loss on ticket price and revenue."
# -> "...have same loss on..." ("has" -> "have").
$ok2 = $d.Content.Find.Execute("runs has same loss", $true, $false, $false, $false, $false,
                                $true, 1, $false, "runs have same loss", 2)
if (-not $ok2) {
    throw "Could not find target text for the 'has' -> 'have' correction."
}
